$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Price Checker Example")
Write-Host $ws1.Name
